$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures to the
# per-job "Leve Profits" sheets, as produced by the scheduled price-sync runner.

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 379.30768
$ws.Range("I9").Value = 198.77777
$ws.Range("K9").Value = 198.77777
$ws.Range("M9").Value = -29.77777
$ws.Range("H53").Value = 887.8333
$ws.Range("J53").Value = 1057.4
$ws.Range("L53").Value = 1057.4
$ws.Range("N53").Value = -2331.4
$ws.Range("H70").Value = 4031.3076
$ws.Range("I70").Value = 2486.75
$ws.Range("K70").Value = 7460.25
$ws.Range("M70").Value = -7190.25
$ws.Range("H73").Value = 4031.3076
$ws.Range("I73").Value = 2486.75
$ws.Range("K73").Value = 7460.25
$ws.Range("M73").Value = -6524.25
$ws.Range("H98").Value = 1759.7693
$ws.Range("J98").Value = 8499
$ws.Range("L98").Value = 8499
$ws.Range("N98").Value = -11495
$ws.Range("H122").Value = 1759.7693
$ws.Range("J122").Value = 8499
$ws.Range("L122").Value = 25497
$ws.Range("N122").Value = -30397

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13184.965
$ws.Range("I32").Value = 15594.061
$ws.Range("J32").Value = 4816.5264
$ws.Range("K32").Value = 15594.061
$ws.Range("L32").Value = 4816.5264
$ws.Range("M32").Value = -15307.061
$ws.Range("N32").Value = -5390.5264
$ws.Range("H80").Value = 150000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 150000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 150000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -151996
$ws.Range("H83").Value = 150000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 150000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 450000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -459984
$ws.Range("H122").Value = 1476.0227
$ws.Range("I122").Value = 1376.7222
$ws.Range("J122").Value = 1922.875
$ws.Range("K122").Value = 4130.1666
$ws.Range("L122").Value = 5768.625
$ws.Range("M122").Value = -1680.1666
$ws.Range("N122").Value = -10668.625
$ws.Range("H132").Value = 41858
$ws.Range("I132").Value = 147700.72
$ws.Range("J132").Value = 6577.095
$ws.Range("K132").Value = 443102.16
$ws.Range("L132").Value = 19731.285
$ws.Range("M132").Value = -440572.16
$ws.Range("N132").Value = -24791.285

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 81071.28999999999
$ws.Range("J135").Value = 81071.28999999999
$ws.Range("L135").Value = 81071.28999999999
$ws.Range("N135").Value = -91211.28999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1472.1177
$ws.Range("I31").Value = 1502.7693
$ws.Range("K31").Value = 1502.7693
$ws.Range("M31").Value = -1207.7693
$ws.Range("H34").Value = 1472.1177
$ws.Range("I34").Value = 1502.7693
$ws.Range("K34").Value = 1502.7693
$ws.Range("M34").Value = -1300.7693
$ws.Range("H58").Value = 129200
$ws.Range("I58").Value = 147255.14
$ws.Range("J58").Value = 2814
$ws.Range("K58").Value = 147255.14
$ws.Range("L58").Value = 2814
$ws.Range("M58").Value = -147052.14
$ws.Range("N58").Value = -3220
$ws.Range("H62").Value = 3857.4443
$ws.Range("I62").Value = 3523.25
$ws.Range("K62").Value = 3523.25
$ws.Range("M62").Value = -2899.25
$ws.Range("H65").Value = 3857.4443
$ws.Range("I65").Value = 3523.25
$ws.Range("K65").Value = 17616.25
$ws.Range("M65").Value = -14496.25
$ws.Range("H105").Value = 1984
$ws.Range("I105").Value = 1995.3334
$ws.Range("J105").Value = 1950
$ws.Range("K105").Value = 1995.3334
$ws.Range("L105").Value = 1950
$ws.Range("M105").Value = -248.3334
$ws.Range("N105").Value = -5444
$ws.Range("H134").Value = 259374.75
$ws.Range("I134").Value = 259374.75
$ws.Range("K134").Value = 778124.25
$ws.Range("M134").Value = -775589.25
$ws.Range("H136").Value = 129200
$ws.Range("I136").Value = 147255.14
$ws.Range("J136").Value = 2814
$ws.Range("K136").Value = 441765.42
$ws.Range("L136").Value = 8442
$ws.Range("M136").Value = -439215.42
$ws.Range("N136").Value = -13542

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 800
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 800
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 2400
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -2746
$ws.Range("H92").Value = 3321.9
$ws.Range("J92").Value = 10244
$ws.Range("L92").Value = 30732
$ws.Range("N92").Value = -33228
$ws.Range("H127").Value = 7881.1665
$ws.Range("J127").Value = 7881.1665
$ws.Range("L127").Value = 23643.4995
$ws.Range("N127").Value = -33563.49950000001
$ws.Range("H132").Value = 2238
$ws.Range("I132").Value = 2441.3333
$ws.Range("J132").Value = 2085.5
$ws.Range("K132").Value = 21971.9997
$ws.Range("L132").Value = 18769.5
$ws.Range("M132").Value = -19441.9997
$ws.Range("N132").Value = -23829.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2663.2778
$ws.Range("I122").Value = 2093.56
$ws.Range("J122").Value = 3958.0908
$ws.Range("K122").Value = 6280.68
$ws.Range("L122").Value = 11874.2724
$ws.Range("M122").Value = -3830.68
$ws.Range("N122").Value = -16774.2724
$ws.Range("H132").Value = 104625
$ws.Range("I132").Value = 116028.89
$ws.Range("K132").Value = 348086.67
$ws.Range("M132").Value = -345556.67

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 73615.7
$ws.Range("I132").Value = 111580.55
$ws.Range("J132").Value = 4013.5
$ws.Range("K132").Value = 334741.65
$ws.Range("L132").Value = 12040.5
$ws.Range("M132").Value = -332211.65
$ws.Range("N132").Value = -17100.5
$ws.Range("H133").Value = 90149.5
$ws.Range("J133").Value = 90149.5
$ws.Range("L133").Value = 90149.5
$ws.Range("N133").Value = -95209.5
$ws.Range("H134").Value = 68997.5
$ws.Range("I134").Value = 68997.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 68997.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -63927.5
$ws.Range("N134").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17874.25
$ws.Range("J41").Value = 17874.25
$ws.Range("L41").Value = 17874.25
$ws.Range("N41").Value = -18654.25
$ws.Range("H45").Value = 41711.57
$ws.Range("J45").Value = 41831
$ws.Range("L45").Value = 41831
$ws.Range("N45").Value = -42813
$ws.Range("H132").Value = 147556
$ws.Range("I132").Value = 204817.4
$ws.Range("K132").Value = 614452.2
$ws.Range("M132").Value = -611922.2
$ws.Range("H136").Value = 4602.4165
$ws.Range("I136").Value = 4248.091
$ws.Range("K136").Value = 12744.273
$ws.Range("M136").Value = -10194.273
$ws.Range("H138").Value = 149995
$ws.Range("I138").Value = 150000
$ws.Range("J138").Value = 149990
$ws.Range("K138").Value = 150000
$ws.Range("L138").Value = 149990
$ws.Range("M138").Value = -144860
$ws.Range("N138").Value = -160270
